# Incidence of HT is (nearly) correct
#
# The "incidence2018_plus" sheet derives monthly incidence from
# prevalence2018 by dividing by a number of years. Rows for ages 35-44
# and 55-120 were divided by 6 and should be divided by 12; rows for
# ages 45-54 were divided by 6 and should be divided by 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("incidence2018_plus")

# Rows 37-46 (ages 35-44): prevalence2018!Cx/6 -> prevalence2018!Cx/12
for ($r = 37; $r -le 46; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/12"
}

# Rows 47-56 (ages 45-54): prevalence2018!Cx/6 -> prevalence2018!Cx/11
for ($r = 47; $r -le 56; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/11"
}

# Rows 57-122 (ages 55-120): prevalence2018!Cx/6 -> prevalence2018!Cx/12
for ($r = 57; $r -le 122; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/12"
}

# Bring the sheet to the front with the same view/selection as the
# author left it (best-effort: scroll so row 23 is at the top, then
# select G41).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G41").Select()
